$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append to the results sheet (script output update).
$data = @(
    @("2025-03-17", "eaux souterraines", 12, 1),
    @("2025-03-17", "eaux souterraines", 14, 2),
    @("2025-03-17", "eaux de surface", 14, 2),
    @("2025-03-17", "ruissellement", 14, 1),
    @("2025-03-17", "développement durable", 16, 1),
    @("2025-03-17", "zone tampon", 65, 1),
    @("2025-03-17", "herbicides", 74, 1),
    @("2025-03-17", "eaux souterraines", 129, 1),
    @("2025-03-17", "eaux souterraines", 130, 1)
)

$r = 13
foreach ($row in $data) {
    # Force column A to be written as text (not auto-parsed into a date
    # serial) by temporarily applying a text number format, then clear
    # the formatting again so the cell ends up styled like its neighbours
    # (no explicit style, just like the rest of the sheet).
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("A$r").ClearFormats()

    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]

    $r++
}
